$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New issue #23 goes into the next empty row (row 19) of the Issues sheet.
# Column A ("Issue ID") holds numeric-looking text ("23"), matching the rest
# of that column which stores issue IDs as text rather than numbers. Using a
# leading apostrophe forces Excel to keep it as text instead of coercing it
# to a number; ClearFormats() then drops the "number stored as text" quote
# -prefix style Excel applies, so the cell keeps the plain default style
# like the rest of the sheet.
$ws.Range("A19").Value = "'23"
$ws.Range("A19").ClearFormats()

$ws.Range("B19").Value = "FR_SALES"
$ws.Range("C19").Value = "open"
$ws.Range("D19").Value = "2025-03-26T06:54:29Z"
$ws.Range("E19").Value = "bug"
